# Refresh the cryptos price/volume snapshot (Coin, Link, Price, Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.842.17"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "3.636.83"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "584.06"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("D6").Value = "175.95"
$ws.Range("E6").Value = "  -3.89%  "
$ws.Range("D7").Value = "3.630.32"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").Value = "6.89"
$ws.Range("E11").Value = "  +16.37%  "
$ws.Range("D12").Value = "0.609"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "48.57"
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "4.225.19"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "676.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.60%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "3.635.20"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "70.888.95"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "17.81"
$ws.Range("E21").Value = "  -4.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "17.16"
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("D25").Value = "100.18"
$ws.Range("E25").Value = "  -4.95%  "
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("D30").Value = "34.72"
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("D31").Value = "9.12"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  -5.13%  "
$ws.Range("D33").Value = "7.56"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("E34").Value = "  -6.59%  "
$ws.Range("D35").Value = "3.99"
$ws.Range("D36").Value = "576.58"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "58.58"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "3.569.55"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").Value = "0.345"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("D45").Value = "34.44"
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("E46").Value = "  -6.30%  "
$ws.Range("E47").Value = "  -4.48%  "
$ws.Range("D48").Value = "2.92"
$ws.Range("E48").Value = "  +3.41%  "
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "137.53"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.95%  "
